$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13 (previously "Box" / "Box position" / "2 limit switch") is repurposed
# to describe the Climber active wheel talon, and its F column value is removed.
$ws.Range("D13").Value = "Climber"
$ws.Range("E13").Value = "Climber active wheel"
$ws.Range("F13").ClearContents()

# Update Climber row: F10 sensor description changes from "Encoder" to "Limit switches"
$ws.Range("F10").Value = "Limit switches"

# Add new row 16 describing the Intake Position talon (CAN id 9)
$ws.Range("A16").Value = "CAN"
$ws.Range("B16").Value = 9
$ws.Range("C16").Value = "Talon SRX"
$ws.Range("D16").Value = "Box"
$ws.Range("E16").Value = "Intake Position"
$ws.Range("F16").Value = "Encoder"

# Restore selection to F10 as in the saved workbook
$ws.Range("F10").Select()
